$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: a cell reference and its new text value. Every value here is
# written as literal text (crypto prices/volume deltas are stored as
# formatted strings in this sheet, e.g. "67.813.12", "  +0.51%  "), so each
# write is prefixed with an apostrophe to force Excel to keep it as text
# instead of auto-converting number-looking values (like "1.00") into a
# numeric cell. The style is reset to "Normal" afterwards so no stray
# quote-prefix / text-format styling is left on the cell.
$updates = @(
    @{ Ref = "D2"; Value = "67.813.12" },
    @{ Ref = "E2"; Value = "  +0.51%  " },
    @{ Ref = "D3"; Value = "2.619.44" },
    @{ Ref = "E3"; Value = "  -0.56%  " },
    @{ Ref = "E4"; Value = "  +0.00%  " },
    @{ Ref = "D5"; Value = "595.89" },
    @{ Ref = "E5"; Value = "  -1.17%  " },
    @{ Ref = "D6"; Value = "152.50" },
    @{ Ref = "E6"; Value = "  -1.02%  " },
    @{ Ref = "E7"; Value = "  +0.01%  " },
    @{ Ref = "E8"; Value = "  -1.61%  " },
    @{ Ref = "D9"; Value = "2.619.29" },
    @{ Ref = "E9"; Value = "  -0.43%  " },
    @{ Ref = "E10"; Value = "  +7.00%  " },
    @{ Ref = "E11"; Value = "  -0.67%  " },
    @{ Ref = "E12"; Value = "  -0.33%  " },
    @{ Ref = "E13"; Value = "  -1.63%  " },
    @{ Ref = "D14"; Value = "27.52" },
    @{ Ref = "E14"; Value = "  -1.68%  " },
    @{ Ref = "D15"; Value = "0.0000188" },
    @{ Ref = "E15"; Value = "  +2.74%  " },
    @{ Ref = "D16"; Value = "3.088.06" },
    @{ Ref = "E16"; Value = "  -0.88%  " },
    @{ Ref = "D17"; Value = "67.733.86" },
    @{ Ref = "E17"; Value = "  +0.33%  " },
    @{ Ref = "D18"; Value = "2.608.48" },
    @{ Ref = "E18"; Value = "  -0.99%  " },
    @{ Ref = "D19"; Value = "371.72" },
    @{ Ref = "E19"; Value = "  +1.88%  " },
    @{ Ref = "D20"; Value = "11.20" },
    @{ Ref = "E20"; Value = "  -0.55%  " },
    @{ Ref = "B21"; Value = "Polkadot" },
    @{ Ref = "C21"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot" },
    @{ Ref = "D21"; Value = "4.23" },
    @{ Ref = "E21"; Value = "  -1.85%  " },
    @{ Ref = "B22"; Value = "Uniswap" },
    @{ Ref = "C22"; Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni" },
    @{ Ref = "D22"; Value = "6.60" },
    @{ Ref = "E22"; Value = "  -13.29%  " },
    @{ Ref = "E23"; Value = "  -3.13%  " },
    @{ Ref = "E24"; Value = "  -4.26%  " },
    @{ Ref = "D25"; Value = "72.91" },
    @{ Ref = "E26"; Value = "  +0.06%  " },
    @{ Ref = "D27"; Value = "9.87" },
    @{ Ref = "E27"; Value = "  -2.08%  " },
    @{ Ref = "D28"; Value = "593.16" },
    @{ Ref = "E28"; Value = "  +1.96%  " },
    @{ Ref = "D29"; Value = "2.747.40" },
    @{ Ref = "E29"; Value = "  -0.81%  " },
    @{ Ref = "E30"; Value = "  -0.44%  " },
    @{ Ref = "E31"; Value = "  -0.10%  " },
    @{ Ref = "E32"; Value = "  -2.33%  " },
    @{ Ref = "D33"; Value = "7.80" },
    @{ Ref = "E33"; Value = "  -1.40%  " },
    @{ Ref = "E34"; Value = "  -0.71%  " },
    @{ Ref = "D35"; Value = "1.00" },
    @{ Ref = "E35"; Value = "  +0.08%  " },
    @{ Ref = "D36"; Value = "0.126" },
    @{ Ref = "E36"; Value = "  -3.22%  " },
    @{ Ref = "E37"; Value = "  -1.59%  " },
    @{ Ref = "D38"; Value = "158.45" },
    @{ Ref = "E38"; Value = "  +0.08%  " },
    @{ Ref = "D39"; Value = "19.15" },
    @{ Ref = "E39"; Value = "  -1.65%  " },
    @{ Ref = "E40"; Value = "  +2.94%  " },
    @{ Ref = "E41"; Value = "  -1.24%  " },
    @{ Ref = "D42"; Value = "5.27" },
    @{ Ref = "E42"; Value = "  -0.70%  " },
    @{ Ref = "E43"; Value = "  +2.75%  " },
    @{ Ref = "D44"; Value = "17.11" },
    @{ Ref = "E44"; Value = "  +4.56%  " },
    @{ Ref = "D45"; Value = "0.999" },
    @{ Ref = "E45"; Value = "  +0.00%  " },
    @{ Ref = "D46"; Value = "40.38" },
    @{ Ref = "E46"; Value = "  -2.05%  " },
    @{ Ref = "B47"; Value = "BabyDogeCoin" },
    @{ Ref = "C47"; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" },
    @{ Ref = "D47"; Value = "0.0₆0299" },
    @{ Ref = "E47"; Value = "  +3.75%  " },
    @{ Ref = "B48"; Value = "Aave" },
    @{ Ref = "C48"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave" },
    @{ Ref = "D48"; Value = "156.36" },
    @{ Ref = "E48"; Value = "  +0.02%  " },
    @{ Ref = "D49"; Value = "3.68" },
    @{ Ref = "E49"; Value = "  -1.54%  " },
    @{ Ref = "E50"; Value = "  -2.35%  " },
    @{ Ref = "E51"; Value = "  -1.72%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}
